# Apply cryptos.xlsx data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '42.615.47'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.288.71'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.70'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.72'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  -3.42%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  -3.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.97'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.24'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.68'
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('D15').Value = '2.643.88'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '2.292.19'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.774'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').Value = '42.509.77'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.76'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').Value = '0.0₃0891'
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.99'
$ws.Range('E21').Value = '  -2.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.82'
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.26'
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.84'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '165.31'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.98'
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.45'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.93'
$ws.Range('E35').Value = '  -3.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.48'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0686'
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('E39').Value = '  -0.94%  '
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E41').Value = '  -2.04%  '
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').Value = '1.990.85'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.97'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.98'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('E47').Value = '  -10.26%  '
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('E49').Value = '  +8.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.41'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').Value = '2.510.81'
$ws.Range('E51').Value = '  -0.46%  '
